$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 22:21"

# Swap country labels where ranking order changed (Santa Lucia now ranks above Macao)
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("A202").Value = "Macao"

# Swap country labels where ranking order changed (Montserrat now ranks above Islas Malvinas)
$ws.Range("A216").Value = "Montserrat"
$ws.Range("A217").Value = "Islas Malvinas"

# Update statistic counts for affected rows
$ws.Range("B4").Value = 8719781
$ws.Range("C4").Value = 55565
$ws.Range("D4").Value = 5685251
$ws.Range("E4").Value = 2805504
$ws.Range("G4").Value = 645
$ws.Range("H4").Value = 229026
$ws.Range("B15").Value = 712412
$ws.Range("C15").Value = 1897
$ws.Range("D15").Value = 643523
$ws.Range("E15").Value = 49998
$ws.Range("G15").Value = 48
$ws.Range("H15").Value = 18891
$ws.Range("B20").Value = 417315
$ws.Range("C20").Value = 13441
$ws.Range("E20").Value = 97025
$ws.Range("G20").Value = 46
$ws.Range("H20").Value = 10090
$ws.Range("B33").Value = 211076
$ws.Range("C33").Value = 1928
$ws.Range("D33").Value = 177450
$ws.Range("E33").Value = 23742
$ws.Range("G33").Value = 22
$ws.Range("H33").Value = 9884
$ws.Range("B51").Value = 101826
$ws.Range("C51").Value = 1210
$ws.Range("D51").Value = 61662
$ws.Range("E51").Value = 38899
$ws.Range("G51").Value = 14
$ws.Range("H51").Value = 1265
$ws.Range("B53").Value = 92229
$ws.Range("C53").Value = 536
$ws.Range("D53").Value = 46118
$ws.Range("E53").Value = 44711
$ws.Range("G53").Value = 4
$ws.Range("H53").Value = 1400
$ws.Range("B94").Value = 20405
$ws.Range("C94").Value = 15
$ws.Range("D94").Value = 20100
$ws.Range("E94").Value = 184
$ws.Range("B105").Value = 12501
$ws.Range("C105").Value = 41
$ws.Range("D105").Value = 10748
$ws.Range("E105").Value = 1620
$ws.Range("B108").Value = 11391
$ws.Range("C108").Value = 33
$ws.Range("D108").Value = 10428
$ws.Range("E108").Value = 926
$ws.Range("B118").Value = 8257
$ws.Range("C118").Value = 15
$ws.Range("D118").Value = 7771
$ws.Range("E118").Value = 250
$ws.Range("B124").Value = 6268
$ws.Range("C124").Value = 133
$ws.Range("D124").Value = 3795
$ws.Range("E124").Value = 2343
$ws.Range("B127").Value = 5831
$ws.Range("C127").Value = 17
$ws.Range("D127").Value = 5485
$ws.Range("B157").Value = 2433
$ws.Range("C157").Value = 19
$ws.Range("D157").Value = 1996
$ws.Range("E157").Value = 372
$ws.Range("B159").Value = 2343
$ws.Range("C159").Value = 3
$ws.Range("D159").Value = 1782
$ws.Range("E159").Value = 487
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 74
$ws.Range("B192").Value = 226
$ws.Range("C192").Value = 2
$ws.Range("E192").Value = 12
$ws.Range("B193").Value = 190
$ws.Range("C193").Value = 2
$ws.Range("E193").Value = 6
$ws.Range("B201").Value = 48
$ws.Range("C201").Value = 6
$ws.Range("D201").Value = 27
$ws.Range("E201").Value = 21
$ws.Range("B202").Value = 46
$ws.Range("D202").Value = 46
$ws.Range("E202").Value = 0
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
$ws.Range("D217").Value = 13
$ws.Range("H217").Value = 0
